$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (headmodel path) was removed for the data rows; only the
# header label "headmodel" in E1 remains. ---
$ws.Range("E2:E9").ClearContents() | Out-Null

# --- Column C (Data path): project folder renamed from "cspAnalysis" to
# "CSPRepo" for every subject row. Processed in row order (split around
# row 12) so the shared-string table ends up in the same order Excel
# would naturally produce. ---
$ws.Range("C2:C11").Replace("cspAnalysis", "CSPRepo") | Out-Null

# Row 12 (subject 020) additionally lost the path separator before the
# filename as part of the same edit.
$ws.Range("C12").Value = "W:\Projects\2018-12 POSTHOCSOURCE Project\analysis_maria\CSPRepo\cleanedsub-020_2.mat"

$ws.Range("C13:C21").Replace("cspAnalysis", "CSPRepo") | Out-Null

# --- Selection moved from C24 to F1 ---
$ws.Range("F1").Select() | Out-Null
